# Generate Report for Handoff
# Replace the old source-file GUID (27cc0499-5efd-4619-b65a-ba3ecb4787eb)
# with the new one (c15e0277-08c4-42d0-9f6c-2f5a303b825f) everywhere it is
# referenced across the three report sheets, refresh the handoff/handback
# timestamps, and keep the Markdown hyperlinks' display text in sync with
# the cell text they annotate.

$wb = $excel.ActiveWorkbook

$oldGuid = "27cc0499-5efd-4619-b65a-ba3ecb4787eb"
$newGuid = "c15e0277-08c4-42d0-9f6c-2f5a303b825f"

$oldZhCnHash = "df2d905a2c74868f15a118899301007199a3cf1c"
$newZhCnHash = "ad6d05f0e6091d8398d032365659974c82ec8922"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/oltest/blob/9ee17aed0146d400e072edb2193db131ab5cf642/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, "", "", "e2e\$newGuid.md")

$wsOverview.Range("G2").Value = "2016-08-13 03:10:25"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")

$wsZhCn.Range("G2").Value = "$newGuid.$newZhCnHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-13 03:10:18"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")

$wsDeDe.Range("G2").Value = "$newGuid.$newZhCnHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-13 03:10:25"
